$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F44").Value = "0.2.0"
$ws.Range("F48").Value = "0.2.0"
$ws.Range("F49").Value = "0.2.0"

$ws.Range("A1:F57").AutoFilter(6, @("0.1.0"), 7)
Write-Host "FilterMode=" $ws.FilterMode
Write-Host "AutoFilterMode=" $ws.AutoFilterMode
